# Bitacora.xlsx - "Modificacion de antecedentes, intro e implementacion"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Cell contents (rows 128-133) ---
$ws.Range("C128").Value = @"
https://prodindu.wordpress.com/revolucion-industrial-4-0/
"@
$ws.Range("E128").Value = @"
png
"@
$ws.Range("D128").Value = @"
cuarta-revoluc
"@
$ws.Range("C129").Value = @"
http://sci-hub.tw/https://doi.org/10.1002/j.1551-8833.2006.tb07609.x
"@
$ws.Range("F129").Value = @"
water trough the years
"@
$ws.Range("G129").Value = @"
pericles_....3398
"@
$ws.Range("G130").Value = @"
'@article{article,
author = {Mays, Larry and Koutsoyiannis, Demetris and Angelakis, A.},
year = {2007},
month = {03},
pages = {},
title = {A brief history of urban water supply in antiquity},
volume = {7},
journal = {Water Science & Technology: Water Supply},
doi = {10.2166/ws.2007.001}
}
"@
$ws.Range("C130").Value = @"
https://www.researchgate.net/publication/228350050_A_brief_history_of_urban_water_supply_in_antiquity
"@
$ws.Range("F130").Value = @"
Historia del agua antigua grecia
"@
$ws.Range("D131").Value = @"
39309-Texto del artículo-48399-2-10-20120628
"@
$ws.Range("E131").Value = @"
pdf
"@
$ws.Range("F131").Value = @"
Aguaa evolucion plantas tratamiento
"@
$ws.Range("H131").Value = @"
Lofrano , G. y Brown , J. (2010). Wastewater management through the ages: A
history of mankind, Science of the Total Environment, 408, 5254 – 5264.
"@
$ws.Range("G131").Value = @"
'@article{article,
author = {Lofrano, Giusy and Brown, Jeanette},
year = {2010},
month = {10},
pages = {5254-64},
title = {Wastewater Management through the Ages: A History of Mankind},
volume = {408},
journal = {The Science of the total environment},
doi = {10.1016/j.scitotenv.2010.07.062}
}
"@
$ws.Range("H132").Value = @"
Cooper , P.F. (2007). Historical aspects of wastewater treatment. In: Lens, P., Seeman,
G., Lettinga, G. (eds). Decentralised sanitation and reuse: concepts, systems
and implementation. IWA Publishing.
"@
$ws.Range("C132").Value = @"
http://www.bvsde.paho.org/bvsacd/leeds/cooper.pdf
"@
$ws.Range("F132").Value = @"
toda la evolucion de agua
"@
$ws.Range("D133").Value = @"
UPS-CT005251
"@
$ws.Range("E133").Value = @"
pdf
"@
$ws.Range("F133").Value = @"
Mariela redes etc
"@

# --- Row heights (grew to fit new wrapped content) ---
$ws.Rows.Item(128).RowHeight = 30.75
$ws.Rows.Item(129).RowHeight = 45.75
$ws.Rows.Item(130).RowHeight = 180.75
$ws.Rows.Item(131).RowHeight = 165.75
$ws.Rows.Item(132).RowHeight = 45.75

# --- View state: scrolled down to the new rows ---
$ws.Activate()
$ws.Range("F134").Select()

